# Weekly update: insert a new day's price record for Mango at
# "Vega Modelo de Temuco", pushing the existing historical rows
# (312..394) down by one (they become 313..395), and populate the
# freshly-inserted row 312 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 312 - this shifts rows
# 312..394 down to 313..395, carrying along their values/styles intact.
$ws.Rows(312).Insert()

# Fill the newly-inserted row 312 with the new record.
# Columns that are constant for every "Vega Modelo de Temuco" / Mango
# row (A, B, C, E, F, G, H, I, J, K, L, Q, T) are restated explicitly
# since this is a brand new row, not a copy.
$ws.Cells.Item(312, 1).Value = 10                                    # Mercado ID
$ws.Cells.Item(312, 2).Value = "Vega Modelo de Temuco"                # Mercado
$ws.Cells.Item(312, 3).Value = "La Araucanía"                         # Región
$ws.Cells.Item(312, 4).Value = 44782                                  # Fecha
$ws.Cells.Item(312, 5).Value = 9                                      # Codreg
$ws.Cells.Item(312, 6).Value = "Fruta"                                # Tipo
$ws.Cells.Item(312, 7).Value = 100108                                 # Producto ID
$ws.Cells.Item(312, 8).Value = "Tropicales y subtropicales"           # Producto
$ws.Cells.Item(312, 9).Value = 100108002                              # Categoría ID
$ws.Cells.Item(312, 10).Value = "Mango"                               # Categoría
$ws.Cells.Item(312, 11).Value = "Sin especificar"                     # Variedad
$ws.Cells.Item(312, 12).Value = "Primera"                             # Calidad
$ws.Cells.Item(312, 13).Value = 125                                   # Volumen
$ws.Cells.Item(312, 14).Value = 10000                                 # Precio mínimo
$ws.Cells.Item(312, 15).Value = 10000                                 # Precio máximo
$ws.Cells.Item(312, 16).Value = 10000                                 # Precio promedio ponderado
$ws.Cells.Item(312, 17).Value = '$/bandeja 4 kilos'                   # Unidad de comercialización
$ws.Cells.Item(312, 18).Value = "Brasil"                              # Origen
$ws.Cells.Item(312, 19).Value = 2500                                  # Precio $/Kg
$ws.Cells.Item(312, 20).Value = 4                                     # Kg / unidad

# Make sure the date cell keeps the same date-time number format used
# by every other "Fecha" cell in this column.
$ws.Cells.Item(312, 4).NumberFormat = $ws.Cells.Item(313, 4).NumberFormat
